$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns (D/E) store numeric- and percent-looking values as
# plain text in the source workbook. Force a Text number format on just those
# cells before assigning so Excel does not reinterpret the strings as numbers
# or percentages (which would change both the stored type and the text).

$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '318.85'
$ws.Range('E2').Value = '4.05%'
$ws.Range('D3').Value = '39.63'
$ws.Range('E3').Value = '2.37%'
$ws.Range('D4').Value = '5.142'
$ws.Range('E4').Value = '1.18%'
$ws.Range('D5').Value = '0.08210'
$ws.Range('E5').Value = '1.96%'
$ws.Range('D6').Value = '2.014'
$ws.Range('E6').Value = '3.79%'
$ws.Range('D7').Value = '8.277'
$ws.Range('E7').Value = '4.32%'
$ws.Range('D8').Value = '4.271'
$ws.Range('E8').Value = '2.00%'
$ws.Range('D9').Value = '0.9341'
$ws.Range('E9').Value = '0.56%'
$ws.Range('D10').Value = '0.1415'
$ws.Range('E10').Value = '-2.95%'
$ws.Range('D11').Value = '0.1993'
$ws.Range('E11').Value = '3.34%'
$ws.Range('D12').Value = '0.09040'
$ws.Range('D13').Value = '0.03582'
$ws.Range('E13').Value = '2.28%'
$ws.Range('D14').Value = '0.09799'
$ws.Range('E14').Value = '-0.02%'
$ws.Range('D15').Value = '0.001392'
$ws.Range('E15').Value = '-0.47%'
$ws.Range('D16').Value = '0.006003'
$ws.Range('E16').Value = '1.65%'
$ws.Range('D17').Value = '3.668'
$ws.Range('E17').Value = '-1.87%'
$ws.Range('E18').Value = '-8.75%'
$ws.Range('D19').Value = '0.3463'
$ws.Range('E19').Value = '-0.01%'
$ws.Range('D20').Value = '0.1277'
$ws.Range('E20').Value = '-2.53%'
$ws.Range('D21').Value = '4.896'
$ws.Range('E21').Value = '2.43%'
$ws.Range('D22').Value = '0.2449'
$ws.Range('E22').Value = '1.93%'
$ws.Range('D23').Value = '0.04328'
$ws.Range('E23').Value = '-0.71%'
$ws.Range('D24').Value = '0.001224'
$ws.Range('E24').Value = '-0.58%'
$ws.Range('D25').Value = '0.004780'
$ws.Range('E25').Value = '11.74%'
$ws.Range('D26').Value = '0.0001299'
$ws.Range('E26').Value = '-0.02%'
$ws.Range('D27').Value = '0.0003996'
$ws.Range('E27').Value = '-10.15%'
$ws.Range('D39').Value = '0.02217'
$ws.Range('E39').Value = '7.14%'
$ws.Range('D40').Value = '0.05264'
$ws.Range('E40').Value = '4.36%'
$ws.Range('D41').Value = '0.007523'
$ws.Range('E41').Value = '0.59%'
$ws.Range('D42').Value = '0.01010'
$ws.Range('E42').Value = '0.08%'
$ws.Range('D43').Value = '0.1377'
$ws.Range('E43').Value = '1.95%'
$ws.Range('D44').Value = '0.002116'
$ws.Range('E44').Value = '-1.05%'
$ws.Range('D45').Value = '0.009863'
$ws.Range('E45').Value = '10.47%'
$ws.Range('D46').Value = '0.00006557'
$ws.Range('E46').Value = '6.02%'
$ws.Range('D47').Value = '0.00000000749'
$ws.Range('E47').Value = '-0.04%'
$ws.Range('B48').Value = 'CoinbaseStockToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D48').Value = '0.001200'
$ws.Range('E48').Value = '-24.92%'
$ws.Range('B49').Value = 'BOLO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D49').Value = '0.002767'
$ws.Range('E49').Value = '-0.83%'
$ws.Range('D50').Value = '0.00002098'
$ws.Range('E50').Value = '-0.04%'
$ws.Range('D51').Value = '0.0001998'
$ws.Range('E51').Value = '-0.04%'
